$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old sample table content, but keep the header row's existing
# (bold/bordered/centered) formatting so it can be reused for the new header.
$ws.Range("A1:D1").ClearContents()
$ws.Range("G6:I9").Clear()

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "EMPLOYEE_ID"
$ws.Range("B1").Value = "EMAIL"
$ws.Range("C1").Value = "FIRST_NAME"
$ws.Range("D1").Value = "LAST_NAME"
$ws.Range("E1").Value = "SALARY"

# Extend the header formatting to the new E column.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Employee IDs -----------------------------------------------------
$ws.Range("A2").Value = 101
$ws.Range("A3").Value = 102
$ws.Range("A4").Value = 103
$ws.Range("A5").Value = 104
$ws.Range("A6").Value = 105

# --- Emails -------------------------------------------------------------
$ws.Range("B2").Value = "john.doe@email.com"
$ws.Range("B3").Value = "jane.smith@email.com"
$ws.Range("B4").Value = "mike.johnson@email.com"
$ws.Range("B5").Value = "emily.white@email.com"
$ws.Range("B6").Value = "david.brown@email.com"

# --- First / last names for rows 3-6 -------------------------------------
$ws.Range("C3").Value = "Jane"
$ws.Range("C4").Value = "Mike"
$ws.Range("C5").Value = "Emily"
$ws.Range("C6").Value = "David"

$ws.Range("D3").Value = "Smith"
$ws.Range("D4").Value = "Johnson"
$ws.Range("D5").Value = "White"
$ws.Range("D6").Value = "Brown"

# --- Salaries -------------------------------------------------------------
$ws.Range("E2").Value = 75000
$ws.Range("E3").Value = 85000.5
$ws.Range("E4").Value = 62000.75
$ws.Range("E5").Value = 92000
$ws.Range("E6").Value = 70000.25

# --- First / last name for row 2 (entered last) ----------------------------
$ws.Range("C2").Value = "Norma"
$ws.Range("D2").Value = "Fisher"

$ws.Range("E13").Select()
